$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.977.46"
$ws.Range("E2").Value = "  -1.56%  "

$ws.Range("D3").Value = "3.333.90"
$ws.Range("E3").Value = "  +2.12%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "576.20"
$ws.Range("E5").Value = "  -1.32%  "

$ws.Range("D6").Value = "183.32"
$ws.Range("E6").Value = "  -0.31%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("E8").Value = "  +0.52%  "

$ws.Range("E9").Value = "  -0.70%  "

$ws.Range("D10").Value = "6.64"
$ws.Range("E10").Value = "  +0.43%  "

$ws.Range("D11").Value = "0.405"
$ws.Range("E11").Value = "  -0.34%  "

$ws.Range("D12").Value = "3.914.64"
$ws.Range("E12").Value = "  +2.21%  "

$ws.Range("E13").Value = "  -0.71%  "

$ws.Range("D14").Value = "27.24"
$ws.Range("E14").Value = "  -0.31%  "

$ws.Range("D15").Value = "67.127.84"
$ws.Range("E15").Value = "  -1.32%  "

$ws.Range("E16").Value = "  -0.34%  "

$ws.Range("D17").Value = "3.334.26"
$ws.Range("E17").Value = "  +2.19%  "

$ws.Range("D18").Value = "444.60"
$ws.Range("E18").Value = "  +6.80%  "

$ws.Range("E19").Value = "  +2.56%  "

$ws.Range("D20").Value = "5.65"
$ws.Range("E20").Value = "  -1.16%  "

$ws.Range("D21").Value = "7.69"
$ws.Range("E21").Value = "  +2.43%  "

$ws.Range("D22").Value = "73.89"
$ws.Range("E22").Value = "  +4.06%  "

$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.10%  "

$ws.Range("D24").Value = "3.483.99"
$ws.Range("E24").Value = "  +2.27%  "

$ws.Range("E25").Value = "  +0.82%  "

$ws.Range("E26").Value = "  +1.98%  "

$ws.Range("D27").Value = "0.195"
$ws.Range("E27").Value = "  +3.67%  "

$ws.Range("D28").Value = "8.98"
$ws.Range("E28").Value = "  -3.27%  "

$ws.Range("E29").Value = "  +0.46%  "

$ws.Range("E30").Value = "  +1.06%  "

$ws.Range("D31").Value = "22.93"
$ws.Range("E31").Value = "  +1.47%  "

$ws.Range("D32").Value = "5.31"
$ws.Range("E32").Value = "  -1.80%  "

$ws.Range("E33").Value = "  -0.06%  "

$ws.Range("D34").Value = "6.78"
$ws.Range("E34").Value = "  -0.70%  "

$ws.Range("E35").Value = "  -0.85%  "

$ws.Range("D36").Value = "161.68"
$ws.Range("E36").Value = "  -1.86%  "

$ws.Range("E37").Value = "  +3.98%  "

$ws.Range("E38").Value = "  +3.11%  "

$ws.Range("E39").Value = "  -1.75%  "

$ws.Range("D40").Value = "2.831.68"
$ws.Range("E40").Value = "  +7.81%  "

$ws.Range("D41").Value = "0.789"
$ws.Range("E41").Value = "  -0.21%  "

$ws.Range("D42").Value = "4.43"
$ws.Range("E42").Value = "  +0.08%  "

$ws.Range("D43").Value = "6.20"
$ws.Range("E43").Value = "  -1.02%  "

$ws.Range("D44").Value = "40.36"

$ws.Range("D45").Value = "0.0670"
$ws.Range("E45").Value = "  -0.41%  "

$ws.Range("D46").Value = "24.38"
$ws.Range("E46").Value = "  +1.10%  "

$ws.Range("E47").Value = "  -2.22%  "

$ws.Range("D48").Value = "321.55"
$ws.Range("E48").Value = "  -3.69%  "

$ws.Range("E49").Value = "  +0.05%  "

$ws.Range("D50").Value = "0.983"
$ws.Range("E50").Value = "  +0.08%  "

$ws.Range("D51").Value = "30.90"
$ws.Range("E51").Value = "  +1.15%  "
